$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-10 22:48:39'
$ws.Range("E3").Value = '2026-02-10 22:48:42'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '95%'
$ws.Range("I3").Value = '26.0 mm'
$ws.Range("L3").Value = '57.2 km/h - 244º 22:06 TU'
$ws.Range("O3").Value = '0.9 °C'
$ws.Range("E4").Value = '2026-02-10 22:48:44'
$ws.Range("E5").Value = '2026-02-10 22:48:46'
$ws.Range("O5").Value = '1.5 °C'
$ws.Range("E6").Value = '2026-02-10 22:48:49'
$ws.Range("E7").Value = '2026-02-10 22:48:51'
$ws.Range("J7").Value = '1004.2 hPa'
$ws.Range("E8").Value = '2026-02-10 22:48:54'
$ws.Range("E9").Value = '2026-02-10 22:48:56'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '96%'
$ws.Range("E10").Value = '2026-02-10 22:48:59'
$ws.Range("O10").Value = '10.5 °C'
$ws.Range("E11").Value = '2026-02-10 22:49:01'
$ws.Range("E12").Value = '2026-02-10 22:49:03'
$ws.Range("E13").Value = '2026-02-10 22:49:06'
$ws.Range("J13").Value = '1006.4 hPa'
$ws.Range("O13").Value = '5.4 °C'
$ws.Range("E14").Value = '2026-02-10 22:49:08'
$ws.Range("O14").Value = '13.1 °C'
$ws.Range("E15").Value = '2026-02-10 22:49:11'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '93%'
$ws.Range("E16").Value = '2026-02-10 22:49:13'
$ws.Range("E17").Value = '2026-02-10 22:49:15'
$ws.Range("I17").Value = '0.4 mm'
$ws.Range("E18").Value = '2026-02-10 22:49:18'
$ws.Range("E19").Value = '2026-02-10 22:49:20'
$ws.Range("O19").Value = '7.0 °C'
$ws.Range("E20").Value = '2026-02-10 22:49:23'
$ws.Range("I20").Value = '12.1 mm'
$ws.Range("E21").Value = '2026-02-10 22:49:25'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '92%'
$ws.Range("I21").Value = '8.9 mm'
$ws.Range("E22").Value = '2026-02-10 22:49:28'
$ws.Range("O22").Value = '-0.4 °C'
$ws.Range("E23").Value = '2026-02-10 22:49:30'
$ws.Range("E24").Value = '2026-02-10 22:49:33'
$ws.Range("O24").Value = '11.1 °C'
$ws.Range("E25").Value = '2026-02-10 22:49:35'
$ws.Range("E26").Value = '2026-02-10 22:49:37'
$ws.Range("J26").Value = '1003.3 hPa'
$ws.Range("E27").Value = '2026-02-10 22:49:40'
$ws.Range("E28").Value = '2026-02-10 22:49:42'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '84%'
$ws.Range("E29").Value = '2026-02-10 22:49:45'
$ws.Range("E30").Value = '2026-02-10 22:49:47'
$ws.Range("E31").Value = '2026-02-10 22:49:50'
$ws.Range("O31").Value = '10.7 °C'
$ws.Range("E32").Value = '2026-02-10 22:49:52'
$ws.Range("O32").Value = '10.8 °C'
$ws.Range("E33").Value = '2026-02-10 22:49:55'
$ws.Range("E34").Value = '2026-02-10 22:49:57'
$ws.Range("E35").Value = '2026-02-10 22:50:00'
$ws.Range("J35").Value = '1005.0 hPa'
$ws.Range("E36").Value = '2026-02-10 22:50:02'
$ws.Range("O36").Value = '10.2 °C'
$ws.Range("E37").Value = '2026-02-10 22:50:05'
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '92%'
$ws.Range("E38").Value = '2026-02-10 22:50:07'
$ws.Range("E39").Value = '2026-02-10 22:50:10'
$ws.Range("L39").Value = '82.1 km/h - 318º 22:05 TU'
$ws.Range("E40").Value = '2026-02-10 22:50:12'
$ws.Range("E41").Value = '2026-02-10 22:50:14'
$ws.Range("J41").Value = '1004.3 hPa'
$ws.Range("L41").Value = '62.3 km/h - 289º 22:13 TU'
$ws.Range("O41").Value = '15.0 °C'
$ws.Range("E42").Value = '2026-02-10 22:50:17'
$ws.Range("E43").Value = '2026-02-10 22:50:19'
$ws.Range("E44").Value = '2026-02-10 22:50:22'
$ws.Range("L44").Value = '65.9 km/h - 204º 22:29 TU'
$ws.Range("E45").Value = '2026-02-10 22:50:24'
$ws.Range("O45").Value = '6.4 °C'
$ws.Range("E46").Value = '2026-02-10 22:50:26'
$ws.Range("O46").Value = '15.2 °C'
